$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M).
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats from column F (the old column D, now shifted) onto the new
# D:E columns so dates / thousands-formatted numbers look right.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = Dec-2018 qtr, E = Sep-2018 qtr) with the
# newly reported figures for every existing data row.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 892300
$ws.Range("E8").Value2 = 1055600
$ws.Range("D9").Value2 = 784300
$ws.Range("E9").Value2 = 911100
$ws.Range("D10").Value2 = 108000
$ws.Range("E10").Value2 = 144500
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 16000
$ws.Range("E14").Value2 = 9300
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 877100
$ws.Range("E17").Value2 = 988200
$ws.Range("D18").Value2 = 15200
$ws.Range("E18").Value2 = 67400
$ws.Range("D20").Value2 = 2700
$ws.Range("E20").Value2 = 4800
$ws.Range("D21").Value2 = 51700
$ws.Range("E21").Value2 = 106500
$ws.Range("D22").Value2 = 4500
$ws.Range("E22").Value2 = 4500
$ws.Range("D23").Value2 = 13400
$ws.Range("E23").Value2 = 67800
$ws.Range("D24").Value2 = 3400
$ws.Range("E24").Value2 = 16300
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 10000
$ws.Range("E26").Value2 = 51500
$ws.Range("D27").Value2 = 6200
$ws.Range("E27").Value2 = 48100
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 400
$ws.Range("E29").Value2 = 7600
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -2700
$ws.Range("E32").Value2 = -4800
$ws.Range("D33").Value2 = 6500
$ws.Range("E33").Value2 = 55700
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 6500
$ws.Range("E35").Value2 = 55700
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 272800
$ws.Range("E41").Value2 = 230300
$ws.Range("D42").Value2 = 312200
$ws.Range("E42").Value2 = 309000
$ws.Range("D43").Value2 = 693000
$ws.Range("E43").Value2 = 832100
$ws.Range("D44").Value2 = 88600
$ws.Range("E44").Value2 = 90800
$ws.Range("D45").Value2 = 48700
$ws.Range("E45").Value2 = 95200
$ws.Range("D46").Value2 = 1415400
$ws.Range("E46").Value2 = 1557300
$ws.Range("D47").Value2 = 120500
$ws.Range("E47").Value2 = 130900
$ws.Range("D48").Value2 = 549700
$ws.Range("E48").Value2 = 560600
$ws.Range("D49").Value2 = 259500
$ws.Range("E49").Value2 = 244700
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 131600
$ws.Range("E52").Value2 = 150300
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 2476600
$ws.Range("E54").Value2 = 2643800
$ws.Range("D57").Value2 = 251500
$ws.Range("E57").Value2 = 316900
$ws.Range("D58").Value2 = 47300
$ws.Range("E58").Value2 = 116800
$ws.Range("D59").Value2 = 379100
$ws.Range("E59").Value2 = 413800
$ws.Range("D60").Value2 = 677800
$ws.Range("E60").Value2 = 847500
$ws.Range("D61").Value2 = 335100
$ws.Range("E61").Value2 = 316900
$ws.Range("D62").Value2 = 66000
$ws.Range("E62").Value2 = 73000
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 1125000
$ws.Range("E66").Value2 = 1282600
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 787400
$ws.Range("E72").Value2 = 786900
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 1351600
$ws.Range("E76").Value2 = 1361300
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 6500
$ws.Range("E81").Value2 = 55700
$ws.Range("D83").Value2 = 33700
$ws.Range("E83").Value2 = 34300
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 71600
$ws.Range("E89").Value2 = 90200
$ws.Range("D91").Value2 = -25000
$ws.Range("E91").Value2 = -49700
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = 41700
$ws.Range("E94").Value2 = -42800
$ws.Range("D96").Value2 = -6100
$ws.Range("E96").Value2 = -5900
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -70600
$ws.Range("E100").Value2 = -12800
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("D102").Value2 = 42800
$ws.Range("E102").Value2 = 34600

# A handful of previously reported figures were also restated in this update.
$ws.Range("H43").Value2 = 583800
$ws.Range("I43").Value2 = 721600
$ws.Range("J43").Value2 = 584100
$ws.Range("H45").Value2 = 36500
$ws.Range("I45").Value2 = 26600
$ws.Range("J45").Value2 = 43600

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Range("D5:M102").EntireColumn.AutoFit() | Out-Null
